# Applies the "Finished Tutorial Battles + Added Dash Pop-up" edits to the
# enemyDatabase worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("enemyDatabase")

# Row 10 (Bunbuku) -----------------------------------------------------
$ws.Range("K10").Value = 0      # EVADE 15 -> 0
$ws.Range("P10").Value = 4      # cryo wk/res 2 -> 4

# Row 11 (Ijiraq) -------------------------------------------------------
$ws.Range("D11").Value = 1      # max_hp 20 -> 1
$ws.Range("F11").Value = 1      # max_stagger 2 -> 1
$ws.Range("I11").Value = 0.6    # SPEED 0.2 -> 0.6
$ws.Range("V11").Value = "lance" # Spell x (root) sword -> lance
$ws.Range("W11").Value = "null"  # Spell x (elem) cryo -> null

# Row 12 (Ijiraq2) -------------------------------------------------------
$ws.Range("D12").Value = 120    # max_hp 65 -> 120
$ws.Range("G12").Value = 1.3    # ATK 1.2 -> 1.3
$ws.Range("I12").Value = 0.4    # SPEED 0.2 -> 0.4
$ws.Range("O12").Value = 3      # agni wk/res 2 -> 3

# Update the active cell selection to match the saved view state.
$ws.Activate()
$ws.Range("K12").Select()
